# Updates the cryptocurrency price/volume table on Sheet1 (rows 2-51)
# to the latest scraped values, per the scheduled GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.976.86'
$ws.Cells.Item(2, 5).Value = '  -0.35%  '
$ws.Cells.Item(3, 4).Value = '2.355.49'
$ws.Cells.Item(3, 5).Value = '  -0.46%  '
$ws.Cells.Item(4, 5).Value = '  +0.06%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '240.04'
$ws.Cells.Item(5, 5).Value = '  -0.31%  '
$ws.Cells.Item(6, 5).Value = '  -2.73%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '74.12'
$ws.Cells.Item(7, 5).Value = '  -2.35%  '
$ws.Cells.Item(8, 5).Value = '  -0.02%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.599'
$ws.Cells.Item(9, 5).Value = '  -3.17%  '
$ws.Cells.Item(10, 5).Value = '  -1.08%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '60.65'
$ws.Cells.Item(11, 5).Value = '  +6.21%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '33.31'
$ws.Cells.Item(12, 5).Value = '  +1.08%  '
$ws.Cells.Item(13, 5).Value = '  +0.53%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '7.27'
$ws.Cells.Item(14, 5).Value = '  -1.68%  '
$ws.Cells.Item(15, 4).Value = '2.704.86'
$ws.Cells.Item(15, 5).Value = '  -0.46%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '16.22'
$ws.Cells.Item(16, 5).Value = '  -2.75%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.908'
$ws.Cells.Item(17, 5).Value = '  -1.19%  '
$ws.Cells.Item(18, 4).Value = '2.347.89'
$ws.Cells.Item(18, 5).Value = '  -0.97%  '
$ws.Cells.Item(19, 4).Value = '43.856.65'
$ws.Cells.Item(19, 5).Value = '  -0.37%  '
$ws.Cells.Item(20, 5).Value = '  -0.49%  '
$ws.Cells.Item(21, 2).Value = 'Litecoin'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '78.18'
$ws.Cells.Item(21, 5).Value = '  +0.77%  '
$ws.Cells.Item(22, 2).Value = 'Uniswap'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '6.64'
$ws.Cells.Item(22, 5).Value = '  -1.04%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '253.04'
$ws.Cells.Item(23, 5).Value = '  -2.24%  '
$ws.Cells.Item(24, 5).Value = '  +0.11%  '
$ws.Cells.Item(25, 5).Value = '  +2.02%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '1.86'
$ws.Cells.Item(26, 5).Value = '  +2.15%  '
$ws.Cells.Item(27, 5).Value = '  -1.22%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '10.46'
$ws.Cells.Item(28, 5).Value = '  -3.57%  '
$ws.Cells.Item(29, 5).Value = '  -2.28%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '175.94'
$ws.Cells.Item(30, 5).Value = '  +0.69%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '22.28'
$ws.Cells.Item(31, 5).Value = '  -3.47%  '
$ws.Cells.Item(32, 5).Value = '  -0.45%  '
$ws.Cells.Item(33, 5).Value = '  -2.60%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.0745'
$ws.Cells.Item(34, 5).Value = '  -2.71%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '5.08'
$ws.Cells.Item(35, 5).Value = '  -5.45%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '5.39'
$ws.Cells.Item(36, 5).Value = '  -0.71%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '3.81'
$ws.Cells.Item(37, 5).Value = '  +0.87%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '2.41'
$ws.Cells.Item(38, 5).Value = '  +0.61%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '6.43'
$ws.Cells.Item(39, 5).Value = '  +0.36%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.0273'
$ws.Cells.Item(40, 5).Value = '  -4.03%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '5.53'
$ws.Cells.Item(41, 5).Value = '  +14.16%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '65.71'
$ws.Cells.Item(42, 5).Value = '  +15.97%  '
$ws.Cells.Item(43, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '19.40'
$ws.Cells.Item(43, 5).Value = '  -1.86%  '
$ws.Cells.Item(44, 2).Value = 'FraxShare'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '9.15'
$ws.Cells.Item(44, 5).Value = '  -1.14%  '
$ws.Cells.Item(45, 2).Value = 'Algorand'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.201'
$ws.Cells.Item(45, 5).Value = '  -4.84%  '
$ws.Cells.Item(46, 2).Value = 'Cronos'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.106'
$ws.Cells.Item(46, 5).Value = '  -3.70%  '
$ws.Cells.Item(47, 5).Value = '  +0.03%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '2.46'
$ws.Cells.Item(48, 5).Value = '  -3.98%  '
$ws.Cells.Item(49, 5).Value = '  -2.29%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '1.15'
$ws.Cells.Item(50, 5).Value = '  -3.24%  '
$ws.Cells.Item(51, 5).Value = '  -3.03%  '
